$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = -0.1190471232299269
$ws.Range("C2").Value = 0.6195113727863492
$ws.Range("D2").Value = 0.6154840829883268
$ws.Range("E2").Value = 0.784527936397632
$ws.Range("F2").Value = 0.7858521091642482

# Row 3 (Q0)
$ws.Range("B3").Value = 0.1110944575514639
$ws.Range("C3").Value = 0.7457144292828872
$ws.Range("D3").Value = 0.8675945666821148
$ws.Range("E3").Value = 0.9314475651812693
$ws.Range("F3").Value = 0.9286440077257753
$ws.Range("G3").Value = 121

# Row 4 (Q1)
$ws.Range("B4").Value = 0.04185042333638055
$ws.Range("C4").Value = 0.6698713821298367
$ws.Range("D4").Value = 0.5416986580756208
$ws.Range("E4").Value = 0.7360018057556794
$ws.Range("F4").Value = 0.7411185024187406
$ws.Range("G4").Value = 59
